$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73, pushing the existing rows 73-75 down to 74-76.
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the new weekly price entry.
$ws.Range("A73").Value2 = 11
$ws.Range("B73").Value2 = "Vega Monumental Concepción"
$ws.Range("C73").Value2 = "Bíobío"
$ws.Range("D73").Value2 = 44706
$ws.Range("E73").Value2 = 8
$ws.Range("F73").Value2 = 100112012
$ws.Range("G73").Value2 = "Espinaca"
$ws.Range("H73").Value2 = "Sin especificar"
$ws.Range("I73").Value2 = "Primera"
$ws.Range("J73").Value2 = 100
$ws.Range("K73").Value2 = 5500
$ws.Range("L73").Value2 = 6000
$ws.Range("M73").Value2 = 5750
$ws.Range("N73").Value2 = "`$/cuna 10 kilos"
$ws.Range("O73").Value2 = "Región Metropolitana"
$ws.Range("P73").Value2 = 575
$ws.Range("Q73").Value2 = 10
$ws.Range("R73").Value2 = "Hortaliza"
